# "foreach si structuri repetitive" - add a new week column of attendance
# (săpt. 11 / column M) for the students who attended, shifting their
# previously-entered attendance mark from column I (săpt. 7, the column
# that had mistakenly been reused) one column to the right so it lines
# up under the correct week, then recording the new week's attendance
# mark of 1 in column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 - Attila Bunta
$ws.Range("I6").Clear()
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 1

# Row 7 - Catalina Madalina Paca
$ws.Range("J7").Clear()
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 2
$ws.Range("M7").Value = 1

# Row 8 - Claudiu Druta
$ws.Range("I8").Clear()
$ws.Range("J8").Value = 2
$ws.Range("M8").Value = 1

# Row 9 - Codrut Avram
$ws.Range("I9").Clear()
$ws.Range("J9").Value = 2

# Row 10 - Daniela Cionca (Marie)
$ws.Range("I10").Clear()
$ws.Range("J10").Value = 2
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 2
$ws.Range("M10").Value = 1

# Row 11 - Delia Negrea
$ws.Range("I11").Clear()
$ws.Range("J11").Value = 2

# Row 13 - Levente Nagy
$ws.Range("I13").Clear()
$ws.Range("J13").Value = 2
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 2
$ws.Range("M13").Value = 1

# Row 14 - Luca Seicaru
$ws.Range("I14").Clear()
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 1

# Row 15 - Mark Pop
$ws.Range("K15").Clear()
$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 1

# Row 19 - Razvan Baroi
$ws.Range("J19").Clear()
$ws.Range("K19").Value = 1

# Row 20 - Silvia Naghi
$ws.Range("J20").Clear()
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 2
$ws.Range("M20").Value = 1

# Row 22 - Victor Lazar
$ws.Range("I22").Clear()
$ws.Range("J22").Value = 2
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 2
$ws.Range("M22").Value = 1

# Update the active selection to reflect where the user ended up (M8)
$ws.Range("M8").Select()
